# Adds the "Keep good code style" section (四、良好的代码风格) content near the
# end of the document, between the existing trailing empty paragraph and the
# paragraph holding the "_GoBack" bookmark, and appends two more blank
# paragraphs after the bookmark paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the empty paragraph that immediately follows the last sentence of
# section 3 ("...不包含尾部。"). That empty paragraph is the anchor right
# before which the new "Keep good code style" material must be inserted.
# ---------------------------------------------------------------------------
$anchor = $d.Content
[void]$anchor.Find.Execute("不包含尾部。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Expand(4)
$anchor.Collapse(0)
$anchor.Expand(4)

# Insert a placeholder empty paragraph right after that anchor paragraph; it
# will be overwritten with the real content in the next step.
$anchor.InsertParagraphAfter()

# Re-locate the same anchor, then step forward onto the freshly inserted
# placeholder paragraph.
$target = $d.Content
[void]$target.Find.Execute("不包含尾部。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Expand(4)
$target.Collapse(0)
$target.Expand(4)
$target.Collapse(0)
$target.Expand(4)

# Replace the placeholder paragraph with the full "Keep good code style"
# block (heading, blank line, bullet paragraphs and the code-comment/html
# ordering examples), each run carrying the eastAsia font hint.
$newSectionXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>四、良好的代码风格</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>适当的空行和缩进。</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>排列整齐的注释：</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">int a = 1;   // </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>注释</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">int b = 11;  // </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>注释</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">int c = 111; // </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>注释</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>语句顺序不能随意，比如与</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> html </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>表单相关联的变量的赋值应该和表单在</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> html </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中的顺序一致。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($newSectionXml)

# ---------------------------------------------------------------------------
# Now append two blank paragraphs right after the paragraph that holds the
# "_GoBack" bookmark (which now immediately follows the block we just
# inserted), before the document's final trailing empty paragraph.
# ---------------------------------------------------------------------------
$bookmarkPara = $d.Content
[void]$bookmarkPara.Find.Execute("中的顺序一致。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPara.Expand(4)
$bookmarkPara.Collapse(0)
$bookmarkPara.Expand(4)

$bookmarkPara.InsertParagraphAfter()
$bookmarkPara.Collapse(0)
$bookmarkPara.Expand(4)

$trailingXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bookmarkPara.InsertXML($trailingXml)
